# Import building IDs as a new "ID" lookup column into CostProgression.xlsx.
# Inserts a new column B (shifting former food/wood/stone/metal columns right),
# labels the sheet's own numeric key column "ID", and fills the new column
# with the building identifiers that correspond to each row's DataID (0-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank column at B; everything from old B onward shifts to C..N.
$ws.Columns("B:B").Insert()

# New header for the inserted column, and relabel the original key column.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "DataID"

# Building identifiers for each data row (A2:A6 = 0..4).
$ws.Range("B2").Value = "building_towncenter"
$ws.Range("B3").Value = "building_farm"
$ws.Range("B4").Value = "building_factory"
$ws.Range("B5").Value = "building_filterationplant"
$ws.Range("B6").Value = "building_house"

# Reflect the user's final selection (the freshly filled-in column).
$ws.Range("B2:B6").Select()
